$wb = $excel.ActiveWorkbook

# Sheet "Hoja1": update the conversion text in A1
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.21 = 12262.2 pesos`n✅ 12262.2 pesos = 3.2 = 970.09 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas": update rate values
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 311.29
$ws2.Range("O10").Value = 3817.1
$ws2.Range("N12").Value = 3830
$ws2.Range("O12").Value = 303
